# Insert a new data row at row 247 (pushing the existing rows 247-272 down
# to 248-273) and populate it with a new "Haba" price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(247).Insert()

$ws.Cells.Item(247, 1).Value  = 6
$ws.Cells.Item(247, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(247, 3).Value  = "Metropolitana"
$ws.Cells.Item(247, 4).Value  = 44769
$ws.Cells.Item(247, 5).Value  = 13
$ws.Cells.Item(247, 6).Value  = 100112026
$ws.Cells.Item(247, 7).Value  = "Haba"
$ws.Cells.Item(247, 8).Value  = "Sin especificar"
$ws.Cells.Item(247, 9).Value  = "Primera"
$ws.Cells.Item(247, 10).Value = 300
$ws.Cells.Item(247, 11).Value = 16000
$ws.Cells.Item(247, 12).Value = 17000
$ws.Cells.Item(247, 13).Value = 16600
$ws.Cells.Item(247, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(247, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(247, 16).Value = 664
$ws.Cells.Item(247, 17).Value = 25
$ws.Cells.Item(247, 18).Value = "Hortaliza"
